# "Generate Report for handoff"
# Updates the localization-status workbook:
#  - Status changes from "Handoff transform failed" to "Ready for handoff"
#    on the Overview sheet and each per-language sheet.
#  - Each per-language sheet gets its "Latest Handoff File" hyperlink cell
#    populated and its "Latest Handoff Datetime" / "Handoff Reason" updated.

$wb = $excel.ActiveWorkbook

$languages = @(
    @{ Sheet = "zh-cn"; File = "7b2d98dc-2cf0-463f-bbcc-d44d90e5e2fa.8058f0369191df01e054fa704f581953ba006c63.zh-cn.xlf"; Datetime = "2016-01-18 04:51:13" },
    @{ Sheet = "de-de"; File = "7b2d98dc-2cf0-463f-bbcc-d44d90e5e2fa.8058f0369191df01e054fa704f581953ba006c63.de-de.xlf"; Datetime = "2016-01-18 04:51:23" }
)

# Overview sheet: the handoff status column for the source file row now
# reports success instead of failure.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Status -> Ready for handoff
    $ws.Range("B2").Value = "Ready for handoff"

    # Latest Handoff File: new hyperlink cell pointing at the generated xlf
    $target = "https://github.com/OpenLocalizationTest/oltest/blob/b2497c52a3503781768181b2e25b733fb97129e8/localization/" + $lang.File
    $ws.Hyperlinks.Add($ws.Range("C2"), $target, "", "", $lang.File)
    # Match the workbook's existing custom hyperlink font (RGB 6495ED) instead
    # of the engine's default themed hyperlink color.
    $ws.Range("C2").Font.Underline = $true
    $ws.Range("C2").Font.Color = 15570276

    # Latest Handoff Datetime: stamp of the handoff that just ran
    $ws.Range("D2").Value = $lang.Datetime

    # Handoff Reason: this file is now included in the handoff
    $ws.Range("H2").Value = "Include"
}
